$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 2 (previously empty inline-string cells, now removed entirely)
$ws.Range("A2:B2").ClearContents()

# Update row 3 content: "nsp" is replaced by "Oxea" in A3, plus new measurement values
$ws.Range("A3").Value = "Oxea"
$ws.Range("B3").Value = 350
$ws.Range("C3").Value = 400

# Move the active selection, matching the saved workbook view state
$ws.Range("D6").Select()
